$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Rebecca's education value (row 3, column C)
$ws.Range("C3").Value = 8

# Remove the "Curtis" row entirely (original row 4); rows below shift up
$ws.Rows.Item(4).Delete()

# Append four new people at rows 9-12
$ws.Range("A9").Value = "Seungyoon"
$ws.Range("B9").Value = "F"
$ws.Range("C9").Value = 22
$ws.Range("D9").Value = "A"
$ws.Range("E9").Value = 44
$ws.Range("F9").Value = $false

$ws.Range("A10").Value = "Jeff"
$ws.Range("B10").Value = "M"
$ws.Range("C10").Value = 22
$ws.Range("D10").Value = "W"
$ws.Range("E10").Value = 45
$ws.Range("F10").Value = $false

$ws.Range("A11").Value = "Joshua"
$ws.Range("B11").Value = "M"
$ws.Range("C11").Value = 22
$ws.Range("D11").Value = "W"
$ws.Range("E11").Value = 37
$ws.Range("F11").Value = $false

$ws.Range("A12").Value = "William"
$ws.Range("B12").Value = "M"
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = "W"
$ws.Range("E12").Value = 11
$ws.Range("F12").Value = $true

$ws.Range("F9:F12").NumberFormat = '"TRUE";"TRUE";"FALSE"'

$ws.Range("A2:A12").Select() | Out-Null
